$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.368.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.56"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.49"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.057.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.234.84"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.610.67"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +4.95%  "

$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.86"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("E35").Value = "  -1.07%  "

$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.833"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.823"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "272.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.944.60"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0223"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.78"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.50%  "
